$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.459.91'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '1.630.77'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''1.000'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').Value = '''304.78'
$ws.Range('E6').Value = '  -1.09%  '
$ws.Range('D7').Value = '''0.3729'
$ws.Range('E7').Value = '  -0.71%  '
$ws.Range('D8').Value = '''0.3652'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('D9').Value = '''51.65'
$ws.Range('E9').Value = '  -2.12%  '
$ws.Range('D10').Value = '''0.08179'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('D11').Value = '''1.226'
$ws.Range('E11').Value = '  -3.98%  '
$ws.Range('D12').Value = '''0.9996'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '''22.55'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('D14').Value = '''6.552'
$ws.Range('E14').Value = '  -1.72%  '
$ws.Range('D15').Value = '''0.00001250'
$ws.Range('E15').Value = '  -2.57%  '
$ws.Range('D16').Value = '''7.273'
$ws.Range('E16').Value = '  -2.27%  '
$ws.Range('D17').Value = '1.631.12'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('D18').Value = '''94.42'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').Value = '''0.06973'
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').Value = '''17.77'
$ws.Range('E20').Value = '  -2.81%  '
$ws.Range('D21').Value = '''6.460'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').Value = '''12.77'
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('D24').Value = '23.467.67'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('D25').Value = '''3.188'
$ws.Range('E25').Value = '  +3.35%  '
$ws.Range('D26').Value = '''2.464'
$ws.Range('E26').Value = '  +1.80%  '
$ws.Range('D27').Value = '''21.45'
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').Value = '''150.45'
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('D29').Value = '''5.337'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = '''134.58'
$ws.Range('E30').Value = '  -1.23%  '
$ws.Range('D31').Value = '1.813.14'
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('D32').Value = '''2.266'
$ws.Range('E32').Value = '  -5.03%  '
$ws.Range('D33').Value = '''6.841'
$ws.Range('E33').Value = '  -0.31%  '
$ws.Range('D34').Value = '''1.023'
$ws.Range('E34').Value = '  +4.52%  '
$ws.Range('D35').Value = '''10.96'
$ws.Range('E35').Value = '  +5.27%  '
$ws.Range('D36').Value = '''0.02784'
$ws.Range('E36').Value = '  -2.54%  '
$ws.Range('D37').Value = '''0.2532'
$ws.Range('E37').Value = '  -0.80%  '
$ws.Range('D38').Value = '''0.08756'
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('D39').Value = '''6.062'
$ws.Range('E39').Value = '  -2.22%  '
$ws.Range('D40').Value = '''0.07134'
$ws.Range('E40').Value = '  -3.55%  '
$ws.Range('D41').Value = '''0.7057'
$ws.Range('E41').Value = '  -1.04%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''1.348'
$ws.Range('E42').Value = '  -2.69%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '''16.28'
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('D44').Value = '''12.33'
$ws.Range('E44').Value = '  -2.04%  '
$ws.Range('D45').Value = '''0.6533'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D46').Value = '''2.335'
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').Value = '''0.9998'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '''3.993'
$ws.Range('E48').Value = '  -1.23%  '
$ws.Range('D49').Value = '''0.08034'
$ws.Range('E49').Value = '  +0.50%  '
$ws.Range('D50').Value = '''1.207'
$ws.Range('E50').Value = '  -0.33%  '
$ws.Range('D51').Value = '''125.07'
$ws.Range('E51').Value = '  -3.78%  '

# Reset style (back to default/no explicit style) for text-forced numeric-looking cells
# to avoid introducing a quotePrefix style that is not part of the original formatting
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
